$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Precio por hora" (price per hour) for each course row: 50 -> 49.5
$ws.Range("J2").Value = 49.5
$ws.Range("J3").Value = 49.5
$ws.Range("J4").Value = 49.5

# Update the active selection to match the author's final cursor position
$ws.Range("I6").Select()
